$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the bottom of the used range (row 8) by shifting row 7
# down - i.e. simply write the final desired 7-row block (rows 2-8) over the
# current 6-row block (rows 2-7), since row 8 previously did not exist.

# Row 2: HLV/Idle at port -> Onshore/Idle at port (E/F cleared)
$ws.Range("B2").Value = "Onshore"
$ws.Range("C2").Value = "Idle at port"
$ws.Range("D2").Value = 25368.5
$ws.Range("E2").Value = $null
$ws.Range("F2").Value = $null

# Row 3: HLV/Transit -> Onshore/None (E/F cleared)
$ws.Range("B3").Value = "Onshore"
$ws.Range("C3").Value = "None"
$ws.Range("D3").Value = 424809
$ws.Range("E3").Value = $null
$ws.Range("F3").Value = $null

# Row 4: Towing Group/Idle at port -> Onshore/Transit (E/F cleared)
$ws.Range("B4").Value = "Onshore"
$ws.Range("C4").Value = "Transit"
$ws.Range("D4").Value = 11088
$ws.Range("E4").Value = $null
$ws.Range("F4").Value = $null

# Row 5: Towing Group/Idle at sea -> Towing Group/Idle at port
$ws.Range("B5").Value = "Towing Group"
$ws.Range("C5").Value = "Idle at port"
$ws.Range("D5").Value = 672579.5244204547
$ws.Range("E5").Value = 0.407310331
$ws.Range("F5").Value = 273948.588715518

# Row 6: Towing Group/Maneuvering -> Towing Group/Idle at sea
$ws.Range("B6").Value = "Towing Group"
$ws.Range("C6").Value = "Idle at sea"
$ws.Range("D6").Value = 108108
$ws.Range("E6").Value = 0.203655165
$ws.Range("F6").Value = 22016.75257782

# Row 7: Towing Group/Transit -> Towing Group/Maneuvering
$ws.Range("B7").Value = "Towing Group"
$ws.Range("C7").Value = "Maneuvering"
$ws.Range("D7").Value = 8316
$ws.Range("E7").Value = 2.698858249
$ws.Range("F7").Value = 22443.705198684

# New Row 8: Towing Group/Transit (brand new row at the end)
# Copy A7 (which carries the bordered/bold/centered style used by column A)
# into A8 so the new row picks up the same formatting, then overwrite value.
$ws.Range("A7").Copy($ws.Range("A8"))
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "Towing Group"
$ws.Range("C8").Value = "Transit"
$ws.Range("D8").Value = 72408.375
$ws.Range("E8").Value = 4.753093345
$ws.Range("F8").Value = 344163.7653347644
$ws.Range("G8").Value = "25 GW (CC)"
